$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.084.85"
$ws.Range("E2").Value = "  +4.80%  "
$ws.Range("D3").Value = "3.238.14"
$ws.Range("E3").Value = "  +2.22%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'578.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.83%  "
$ws.Range("D6").Value = "'176.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.33%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("D9").Value = "3.233.70"
$ws.Range("E9").Value = "  +2.15%  "
$ws.Range("E10").Value = "  +4.76%  "
$ws.Range("E11").Value = "  +1.51%  "
$ws.Range("D12").Value = "'0.406"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.06%  "
$ws.Range("D13").Value = "3.800.35"
$ws.Range("E13").Value = "  +2.24%  "
$ws.Range("E14").Value = "  +1.99%  "
$ws.Range("D15").Value = "'27.79"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.50%  "
$ws.Range("D16").Value = "67.017.75"
$ws.Range("E16").Value = "  +4.48%  "
$ws.Range("D17").Value = "'0.0000167"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.81%  "
$ws.Range("D18").Value = "3.238.97"
$ws.Range("E18").Value = "  +2.32%  "
$ws.Range("E19").Value = "  +2.88%  "
$ws.Range("D20").Value = "'13.25"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.29%  "
$ws.Range("D21").Value = "'367.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.47%  "
$ws.Range("E22").Value = "  +3.57%  "
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").Value = "'70.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.95%  "
$ws.Range("D25").Value = "3.381.51"
$ws.Range("E25").Value = "  +2.22%  "
$ws.Range("E26").Value = "  +0.92%  "
$ws.Range("E27").Value = "  +0.79%  "
$ws.Range("D28").Value = "'9.70"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.67%  "
$ws.Range("E29").Value = "  +1.85%  "
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("E31").Value = "  +4.74%  "
$ws.Range("D32").Value = "'5.58"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.46%  "
$ws.Range("D33").Value = "'22.40"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.62%  "
$ws.Range("E34").Value = "  -0.11%  "
$ws.Range("D35").Value = "'173.69"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +10.46%  "
$ws.Range("D36").Value = "'1.22"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.36%  "
$ws.Range("D37").Value = "'6.74"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.45%  "
$ws.Range("E38").Value = "  +5.20%  "
$ws.Range("D39").Value = "'0.848"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.74%  "
$ws.Range("E40").Value = "  +10.51%  "
$ws.Range("D41").Value = "'26.56"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Value = "'2.58"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.47%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "'6.38"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.23%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "2.708.90"
$ws.Range("E44").Value = "  +2.79%  "
$ws.Range("D45").Value = "'4.28"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.60%  "
$ws.Range("D46").Value = "'40.37"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.74%  "
$ws.Range("D47").Value = "'0.0668"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.64%  "
$ws.Range("E48").Value = "  +5.13%  "
$ws.Range("D49").Value = "'333.86"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.28%  "
$ws.Range("E50").Value = "  +3.35%  "
$ws.Range("E51").Value = "  +2.35%  "
